$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.430.38'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '1.583.60'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'213.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = "'0.491"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = "'44.52"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("D9").Value = "'23.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("D10").Value = "'0.247"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("E11").Value = '  -1.76%  '
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("D13").Value = '1.810.49'
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").Value = '1.582.93'
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").Value = "'3.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("D17").Value = '28.435.66'
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("D18").Value = "'62.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.60%  '
$ws.Range("D19").Value = "'229.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("D20").Value = "'7.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("E21").Value = '  -2.26%  '
$ws.Range("E23").Value = '  -3.20%  '
$ws.Range("D24").Value = "'9.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  +2.72%  '
$ws.Range("D26").Value = "'152.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").Value = "'15.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.17%  '
$ws.Range("D28").Value = "'6.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.78%  '
$ws.Range("E29").Value = '  -1.84%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").Value = "'0.0481"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.46%  '
$ws.Range("E32").Value = '  -1.17%  '
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("E34").Value = '  -2.54%  '
$ws.Range("D35").Value = '1.395.92'
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("E36").Value = '  +7.46%  '
$ws.Range("D37").Value = "'1.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.17%  '
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("E39").Value = '  +0.72%  '
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("E41").Value = '  -3.25%  '
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").Value = "'0.791"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("E45").Value = '  -3.39%  '
$ws.Range("E46").Value = '  -1.45%  '
$ws.Range("D47").Value = "'0.957"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("D48").Value = "'62.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("D49").Value = '1.721.88'
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").Value = "'86.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("E51").Value = '  -1.38%  '
